$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 61
$ws.Range("B61").Value = 6905571
$ws.Range("C61").Value = 'Peru Liga 1'
$ws.Range("D61").Value = 45130.72916666666
$ws.Range("E61").Value = 'FBC Melgar'
$ws.Range("F61").Value = 'Sporting Cristal'
$ws.Range("G61").Value = 1
$ws.Range("H61").Value = 1
$ws.Range("I61").Value = 1
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 'D'
$ws.Range("L61").Value = 2.1
$ws.Range("M61").Value = 3.4
$ws.Range("N61").Value = 3
$ws.Range("O61").Value = 1.75
$ws.Range("P61").Value = 3.8
$ws.Range("Q61").Value = 4.75
$ws.Range("R61").Value = -0.75
$ws.Range("S61").Value = 1.95
$ws.Range("T61").Value = 1.85
$ws.Range("U61").Value = 2.5
$ws.Range("V61").Value = 1.95
$ws.Range("W61").Value = 1.85
$ws.Range("X61").Value = -1
$ws.Range("Y61").Value = 2.8
$ws.Range("Z61").Value = -1
$ws.Range("AA61").Value = -1
$ws.Range("AB61").Value = 0.8500000000000001
$ws.Range("AC61").Value = -1
$ws.Range("AD61").Value = 0.8500000000000001

# Row 62
$ws.Range("B62").Value = 6905578
$ws.Range("C62").Value = 'Peru Liga 1'
$ws.Range("D62").Value = 45130.72916666666
$ws.Range("E62").Value = 'AD Tarma'
$ws.Range("F62").Value = 'Atletico Grau'
$ws.Range("G62").Value = 1
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 1
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 'H'
$ws.Range("L62").Value = 1.75
$ws.Range("M62").Value = 3.6
$ws.Range("N62").Value = 4
$ws.Range("O62").Value = 1.571
$ws.Range("P62").Value = 4.2
$ws.Range("Q62").Value = 5.75
$ws.Range("R62").Value = -1
$ws.Range("S62").Value = 1.975
$ws.Range("T62").Value = 1.825
$ws.Range("U62").Value = 2.5
$ws.Range("V62").Value = 1.8
$ws.Range("W62").Value = 2
$ws.Range("X62").Value = 0.571
$ws.Range("Y62").Value = -1
$ws.Range("Z62").Value = -1
$ws.Range("AA62").Value = 0
$ws.Range("AB62").Value = 0
$ws.Range("AC62").Value = -1
$ws.Range("AD62").Value = 1

# Row 156
$ws.Range("B156").Value = 7211641
$ws.Range("C156").Value = 'Peru Liga 1'
$ws.Range("D156").Value = 45198.70833333334
$ws.Range("E156").Value = 'Sport Huancayo'
$ws.Range("F156").Value = 'Deportivo Municipal'
$ws.Range("G156").Value = 2
$ws.Range("H156").Value = 0
$ws.Range("I156").Value = 0
$ws.Range("J156").Value = 0
$ws.Range("K156").Value = 'H'
$ws.Range("L156").Value = 1.125
$ws.Range("M156").Value = 7
$ws.Range("N156").Value = 17
$ws.Range("O156").Value = 1.166
$ws.Range("P156").Value = 6.5
$ws.Range("Q156").Value = 12
$ws.Range("R156").Value = -2
$ws.Range("S156").Value = 1.775
$ws.Range("T156").Value = 2.025
$ws.Range("U156").Value = 3.5
$ws.Range("V156").Value = 1.9
$ws.Range("W156").Value = 1.9
$ws.Range("X156").Value = 0.1659999999999999
$ws.Range("Y156").Value = -1
$ws.Range("Z156").Value = -1
$ws.Range("AA156").Value = 0
$ws.Range("AB156").Value = 0
$ws.Range("AC156").Value = -1
$ws.Range("AD156").Value = 0.8999999999999999

# Row 157
$ws.Range("B157").Value = 7211640
$ws.Range("C157").Value = 'Peru Liga 1'
$ws.Range("D157").Value = 45198.70833333334
$ws.Range("E157").Value = 'UTC Cajamarca'
$ws.Range("F157").Value = 'Sport Boys'
$ws.Range("G157").Value = 1
$ws.Range("H157").Value = 1
$ws.Range("I157").Value = 0
$ws.Range("J157").Value = 0
$ws.Range("K157").Value = 'D'
$ws.Range("L157").Value = 1.615
$ws.Range("M157").Value = 3.75
$ws.Range("N157").Value = 5
$ws.Range("O157").Value = 1.5
$ws.Range("P157").Value = 4.2
$ws.Range("Q157").Value = 6.5
$ws.Range("R157").Value = -1
$ws.Range("S157").Value = 1.8
$ws.Range("T157").Value = 2.05
$ws.Range("U157").Value = 2.5
$ws.Range("V157").Value = 1.875
$ws.Range("W157").Value = 1.975
$ws.Range("X157").Value = -1
$ws.Range("Y157").Value = 3.2
$ws.Range("Z157").Value = -1
$ws.Range("AA157").Value = -1
$ws.Range("AB157").Value = 1.05
$ws.Range("AC157").Value = -1
$ws.Range("AD157").Value = 0.9750000000000001

# Row 184
$ws.Range("B184").Value = 7384626
$ws.Range("C184").Value = 'Peru Liga 1'
$ws.Range("D184").Value = 45228.70833333334
$ws.Range("E184").Value = 'Sporting Cristal'
$ws.Range("F184").Value = 'Alianza Atletico'
$ws.Range("G184").Value = 3
$ws.Range("H184").Value = 0
$ws.Range("I184").Value = 3
$ws.Range("J184").Value = 0
$ws.Range("K184").Value = 'H'
$ws.Range("L184").Value = 1.3
$ws.Range("M184").Value = 5
$ws.Range("N184").Value = 9
$ws.Range("O184").Value = 1.166
$ws.Range("P184").Value = 6.5
$ws.Range("Q184").Value = 13
$ws.Range("R184").Value = -2
$ws.Range("S184").Value = 1.85
$ws.Range("T184").Value = 1.95
$ws.Range("U184").Value = 3.25
$ws.Range("V184").Value = 2
$ws.Range("W184").Value = 1.8
$ws.Range("X184").Value = 0.1659999999999999
$ws.Range("Y184").Value = -1
$ws.Range("Z184").Value = -1
$ws.Range("AA184").Value = 0.8500000000000001
$ws.Range("AB184").Value = -1
$ws.Range("AC184").Value = -0.5
$ws.Range("AD184").Value = 0.4

# Row 185
$ws.Range("B185").Value = 7384627
$ws.Range("C185").Value = 'Peru Liga 1'
$ws.Range("D185").Value = 45228.70833333334
$ws.Range("E185").Value = 'Universitario de Deportes'
$ws.Range("F185").Value = 'Sport Huancayo'
$ws.Range("G185").Value = 2
$ws.Range("H185").Value = 0
$ws.Range("I185").Value = 1
$ws.Range("J185").Value = 0
$ws.Range("K185").Value = 'H'
$ws.Range("L185").Value = 1.25
$ws.Range("M185").Value = 5
$ws.Range("N185").Value = 12
$ws.Range("O185").Value = 1.181
$ws.Range("P185").Value = 6
$ws.Range("Q185").Value = 13
$ws.Range("R185").Value = -1.75
$ws.Range("S185").Value = 1.8
$ws.Range("T185").Value = 2
$ws.Range("U185").Value = 2.75
$ws.Range("V185").Value = 1.85
$ws.Range("W185").Value = 1.95
$ws.Range("X185").Value = 0.181
$ws.Range("Y185").Value = -1
$ws.Range("Z185").Value = -1
$ws.Range("AA185").Value = 0.4
$ws.Range("AB185").Value = -0.5
$ws.Range("AC185").Value = -1
$ws.Range("AD185").Value = 0.95

# Row 186
$ws.Range("B186").Value = 7384629
$ws.Range("C186").Value = 'Peru Liga 1'
$ws.Range("D186").Value = 45228.70833333334
$ws.Range("E186").Value = 'Deportivo Garcilaso'
$ws.Range("F186").Value = 'Alianza Lima'
$ws.Range("G186").Value = 0
$ws.Range("H186").Value = 1
$ws.Range("I186").Value = 0
$ws.Range("J186").Value = 1
$ws.Range("K186").Value = 'A'
$ws.Range("L186").Value = 2.625
$ws.Range("M186").Value = 3.3
$ws.Range("N186").Value = 2.5
$ws.Range("O186").Value = 2.7
$ws.Range("P186").Value = 3.4
$ws.Range("Q186").Value = 2.375
$ws.Range("R186").Value = 0
$ws.Range("S186").Value = 2.025
$ws.Range("T186").Value = 1.775
$ws.Range("U186").Value = 2.25
$ws.Range("V186").Value = 1.825
$ws.Range("W186").Value = 1.975
$ws.Range("X186").Value = -1
$ws.Range("Y186").Value = -1
$ws.Range("Z186").Value = 1.375
$ws.Range("AA186").Value = -1
$ws.Range("AB186").Value = 0.7749999999999999
$ws.Range("AC186").Value = -1
$ws.Range("AD186").Value = 0.9750000000000001

# Row 187
$ws.Range("B187").Value = 7384628
$ws.Range("C187").Value = 'Peru Liga 1'
$ws.Range("D187").Value = 45228.70833333334
$ws.Range("E187").Value = 'Deportivo Binacional'
$ws.Range("F187").Value = 'FBC Melgar'
$ws.Range("G187").Value = 1
$ws.Range("H187").Value = 2
$ws.Range("I187").Value = 1
$ws.Range("J187").Value = 1
$ws.Range("K187").Value = 'A'
$ws.Range("L187").Value = 2.75
$ws.Range("M187").Value = 3.3
$ws.Range("N187").Value = 2.375
$ws.Range("O187").Value = 3.3
$ws.Range("P187").Value = 3.6
$ws.Range("Q187").Value = 2
$ws.Range("R187").Value = 0.5
$ws.Range("S187").Value = 1.8
$ws.Range("T187").Value = 2
$ws.Range("U187").Value = 2.75
$ws.Range("V187").Value = 1.975
$ws.Range("W187").Value = 1.875
$ws.Range("X187").Value = -1
$ws.Range("Y187").Value = -1
$ws.Range("Z187").Value = 1
$ws.Range("AA187").Value = -1
$ws.Range("AB187").Value = 1
$ws.Range("AC187").Value = 0.4875
$ws.Range("AD187").Value = -0.5

# Row 188
$ws.Range("B188").Value = 7384625
$ws.Range("C188").Value = 'Peru Liga 1'
$ws.Range("D188").Value = 45228.70833333334
$ws.Range("E188").Value = 'AD Tarma'
$ws.Range("F188").Value = 'Carlos Manucci'
$ws.Range("G188").Value = 0
$ws.Range("H188").Value = 0
$ws.Range("I188").Value = 0
$ws.Range("J188").Value = 0
$ws.Range("K188").Value = 'D'
$ws.Range("L188").Value = 1.5
$ws.Range("M188").Value = 3.75
$ws.Range("N188").Value = 7
$ws.Range("O188").Value = 1.363
$ws.Range("P188").Value = 4.333
$ws.Range("Q188").Value = 9.5
$ws.Range("R188").Value = -1.25
$ws.Range("S188").Value = 1.875
$ws.Range("T188").Value = 1.925
$ws.Range("U188").Value = 2.5
$ws.Range("V188").Value = 1.8
$ws.Range("W188").Value = 2
$ws.Range("X188").Value = -1
$ws.Range("Y188").Value = 3.333
$ws.Range("Z188").Value = -1
$ws.Range("AA188").Value = -1
$ws.Range("AB188").Value = 0.925
$ws.Range("AC188").Value = -1
$ws.Range("AD188").Value = 1

# Row 228
$ws.Range("B228").Value = 7818816
$ws.Range("C228").Value = 'Peru Liga 1'
$ws.Range("D228").Value = 45346.70833333334
$ws.Range("E228").Value = 'UTC Cajamarca'
$ws.Range("F228").Value = 'Universitario de Deportes'
$ws.Range("G228").Value = 0
$ws.Range("H228").Value = 0
$ws.Range("I228").Value = 0
$ws.Range("J228").Value = 0
$ws.Range("K228").Value = 'D'
$ws.Range("L228").Value = 3.3
$ws.Range("M228").Value = 3.3
$ws.Range("N228").Value = 2.1
$ws.Range("O228").Value = 4.5
$ws.Range("P228").Value = 3.2
$ws.Range("Q228").Value = 1.95
$ws.Range("R228").Value = 0.5
$ws.Range("S228").Value = 2
$ws.Range("T228").Value = 1.85
$ws.Range("U228").Value = 2
$ws.Range("V228").Value = 1.775
$ws.Range("W228").Value = 2.1
$ws.Range("X228").Value = -1
$ws.Range("Y228").Value = 2.2
$ws.Range("Z228").Value = -1
$ws.Range("AA228").Value = 1
$ws.Range("AB228").Value = -1
$ws.Range("AC228").Value = -1
$ws.Range("AD228").Value = 1.1

# Row 229
$ws.Range("B229").Value = 7818817
$ws.Range("C229").Value = 'Peru Liga 1'
$ws.Range("D229").Value = 45346.70833333334
$ws.Range("E229").Value = 'Sport Boys'
$ws.Range("F229").Value = 'Cusco FC'
$ws.Range("G229").Value = 3
$ws.Range("H229").Value = 0
$ws.Range("I229").Value = 2
$ws.Range("J229").Value = 0
$ws.Range("K229").Value = 'H'
$ws.Range("L229").Value = 2.2
$ws.Range("M229").Value = 3.2
$ws.Range("N229").Value = 3.2
$ws.Range("O229").Value = 1.6
$ws.Range("P229").Value = 3.75
$ws.Range("Q229").Value = 5.75
$ws.Range("R229").Value = -0.75
$ws.Range("S229").Value = 1.85
$ws.Range("T229").Value = 2
$ws.Range("U229").Value = 2.5
$ws.Range("V229").Value = 1.975
$ws.Range("W229").Value = 1.875
$ws.Range("X229").Value = 0.6000000000000001
$ws.Range("Y229").Value = -1
$ws.Range("Z229").Value = -1
$ws.Range("AA229").Value = 0.8500000000000001
$ws.Range("AB229").Value = -1
$ws.Range("AC229").Value = 0.9750000000000001
$ws.Range("AD229").Value = -1

# Row 252
$ws.Range("B252").Value = 7883367
$ws.Range("C252").Value = 'Peru Liga 1'
$ws.Range("D252").Value = 45361.70833333334
$ws.Range("E252").Value = 'Sport Boys'
$ws.Range("F252").Value = 'Cesar Vallejo'
$ws.Range("G252").Value = 2
$ws.Range("H252").Value = 0
$ws.Range("I252").Value = 2
$ws.Range("J252").Value = 0
$ws.Range("K252").Value = 'H'
$ws.Range("L252").Value = 2.2
$ws.Range("M252").Value = 3.3
$ws.Range("N252").Value = 3.1
$ws.Range("O252").Value = 2.4
$ws.Range("P252").Value = 3.25
$ws.Range("Q252").Value = 2.625
$ws.Range("R252").Value = 0
$ws.Range("S252").Value = 1.8
$ws.Range("T252").Value = 2
$ws.Range("U252").Value = 2.5
$ws.Range("V252").Value = 2.025
$ws.Range("W252").Value = 1.825
$ws.Range("X252").Value = 1.4
$ws.Range("Y252").Value = -1
$ws.Range("Z252").Value = -1
$ws.Range("AA252").Value = 0.8
$ws.Range("AB252").Value = -1
$ws.Range("AC252").Value = -1
$ws.Range("AD252").Value = 0.825

# Row 253
$ws.Range("B253").Value = 7882752
$ws.Range("C253").Value = 'Peru Liga 1'
$ws.Range("D253").Value = 45361.70833333334
$ws.Range("E253").Value = 'Sport Huancayo'
$ws.Range("F253").Value = 'Union Comercio'
$ws.Range("G253").Value = 2
$ws.Range("H253").Value = 2
$ws.Range("I253").Value = 0
$ws.Range("J253").Value = 1
$ws.Range("K253").Value = 'D'
$ws.Range("L253").Value = 1.3
$ws.Range("M253").Value = 4.5
$ws.Range("N253").Value = 10
$ws.Range("O253").Value = 1.3
$ws.Range("P253").Value = 4.333
$ws.Range("Q253").Value = 9.5
$ws.Range("R253").Value = -1.5
$ws.Range("S253").Value = 1.95
$ws.Range("T253").Value = 1.85
$ws.Range("U253").Value = 2.75
$ws.Range("V253").Value = 1.975
$ws.Range("W253").Value = 1.825
$ws.Range("X253").Value = -1
$ws.Range("Y253").Value = 3.333
$ws.Range("Z253").Value = -1
$ws.Range("AA253").Value = -1
$ws.Range("AB253").Value = 0.8500000000000001
$ws.Range("AC253").Value = 0.9750000000000001
$ws.Range("AD253").Value = -1

# Row 305
$ws.Range("B305").Value = 8132794
$ws.Range("C305").Value = 'Peru Liga 1'
$ws.Range("D305").Value = 45410.91666666666
$ws.Range("E305").Value = 'Universitario de Deportes'
$ws.Range("F305").Value = 'Comerciantes Unidos'
$ws.Range("G305").Value = 6
$ws.Range("H305").Value = 0
$ws.Range("I305").Value = 4
$ws.Range("J305").Value = 0
$ws.Range("K305").Value = 'H'
$ws.Range("L305").Value = 1.166
$ws.Range("M305").Value = 7
$ws.Range("N305").Value = 12
$ws.Range("O305").Value = 1.2
$ws.Range("P305").Value = 6.5
$ws.Range("Q305").Value = 13
$ws.Range("R305").Value = -1.75
$ws.Range("S305").Value = 1.8
$ws.Range("T305").Value = 2
$ws.Range("U305").Value = 3
$ws.Range("V305").Value = 1.975
$ws.Range("W305").Value = 1.825
$ws.Range("X305").Value = 0.2
$ws.Range("Y305").Value = -1
$ws.Range("Z305").Value = -1
$ws.Range("AA305").Value = 0.8
$ws.Range("AB305").Value = -1
$ws.Range("AC305").Value = 0.9750000000000001
$ws.Range("AD305").Value = -1

# Row 306
$ws.Range("B306").Value = 8042217
$ws.Range("C306").Value = 'Peru Liga 1'
$ws.Range("D306").Value = 45410.91666666666
$ws.Range("E306").Value = 'Deportivo Garcilaso'
$ws.Range("F306").Value = 'CD Los Chankas'
$ws.Range("G306").Value = 1
$ws.Range("H306").Value = 1
$ws.Range("I306").Value = 0
$ws.Range("J306").Value = 1
$ws.Range("K306").Value = 'D'
$ws.Range("L306").Value = 2
$ws.Range("M306").Value = 3.3
$ws.Range("N306").Value = 3.75
$ws.Range("O306").Value = 1.727
$ws.Range("P306").Value = 3.6
$ws.Range("Q306").Value = 4.5
$ws.Range("R306").Value = -0.75
$ws.Range("S306").Value = 1.9
$ws.Range("T306").Value = 1.9
$ws.Range("U306").Value = 2.75
$ws.Range("V306").Value = 1.8
$ws.Range("W306").Value = 2
$ws.Range("X306").Value = -1
$ws.Range("Y306").Value = 2.6
$ws.Range("Z306").Value = -1
$ws.Range("AA306").Value = -1
$ws.Range("AB306").Value = 0.8999999999999999
$ws.Range("AC306").Value = -1
$ws.Range("AD306").Value = 1
